# Fruta / hortaliza, semanal
# Pimiento sheet: a new weekly price report is inserted for
# "Zafiro rojo" / "Zafiro verde" (Región de Arica y Parinacota), while the
# two most-recent existing rows (295/296) are updated in place to the
# latest week's figures. Net effect: 2 new rows appear, the sheet grows
# from A1:R318 to A1:R320.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows right after the current last "Zafiro"
# entries (rows 295-296), pushing everything from the old row 297 onward
# down by two rows.
$ws.Range("A297:A298").EntireRow.Insert()

# The two new rows (297-298) start life as exact copies of the rows they
# are being inserted after (295-296) - i.e. the data that is about to be
# superseded by this week's update.
$ws.Range("A295:R295").Copy($ws.Range("A297:R297"))
$ws.Range("A296:R296").Copy($ws.Range("A298:R298"))

# Now refresh rows 295-296 with this week's figures.
$ws.Range("D295").Value = 44461
$ws.Range("H295").Value = "Zafiro rojo"
$ws.Range("J295").Value = 400
$ws.Range("K295").Value = 35000
$ws.Range("L295").Value = 35000
$ws.Range("M295").Value = 35000
$ws.Range("O295").Value = "Región de Arica y Parinacota"
$ws.Range("P295").Value = 2333

$ws.Range("D296").Value = 44461
$ws.Range("H296").Value = "Zafiro verde"
$ws.Range("J296").Value = 300
$ws.Range("K296").Value = 28000
$ws.Range("L296").Value = 28000
$ws.Range("M296").Value = 28000
$ws.Range("O296").Value = "Región de Arica y Parinacota"
$ws.Range("P296").Value = 1867
